$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G18").Value = 44989
$ws.Range("G18").NumberFormat = "m/d/yyyy"
